# "Generate Report for Archive" — refresh the localization-status report:
# the zh-cn / de-de handoff status moves from "Ready for handoff" to
# "In Translation", and the Status columns that held the old (longer)
# text are re-fitted to the new (shorter) text's width.

$wb = $excel.ActiveWorkbook

$newStatus = "In Translation"

# --- Overview sheet: one status cell per locale column (E, F) ---
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2").Value = $newStatus
$overview.Range("F2").Value = $newStatus
$overview.Columns("E").ColumnWidth = 12.5
$overview.Columns("F").ColumnWidth = 12.5

# --- zh-cn sheet: Status column (C) ---
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C2").Value = $newStatus
$zhcn.Columns("C").ColumnWidth = 12.5

# --- de-de sheet: Status column (C) ---
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C2").Value = $newStatus
$dede.Columns("C").ColumnWidth = 12.5
